$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Series Cast")

# Rows 53 and 54: swap Name/Screen Name between the two rows
$ws.Range("A53").Value = "Megan Heffernan"
$ws.Range("B53").Value = "Breast Milk Woman"
$ws.Range("A54").Value = "Nicky B"
$ws.Range("B54").Value = "Fat Woman"

# Rows 90-92: rotate the three rows' Name/Screen Name values
$ws.Range("A90").Value = "Ken Byrd"
$ws.Range("B90").Value = "Passerby (uncredited)"
$ws.Range("A91").Value = "Myra Ford"
$ws.Range("B91").Value = "Checkout Assistant (uncredited)"
$ws.Range("A92").Value = "Karen-J Sear"
$ws.Range("B92").Value = "Passer by (uncredited)"

# Rows 97 and 98: swap Name/Screen Name between the two rows
$ws.Range("A97").Value = "Gleb Smatko"
$ws.Range("B97").Value = "Cafe Patron (uncredited)"
$ws.Range("A98").Value = "Marta Glowacka-Escote"
$ws.Range("B98").Value = "Coffee shop guest (uncredited)"
